$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "243.52"

# Row 4
$ws.Range("B4").Value = "LEO"
$ws.Range("C4").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "3.577"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "3LEOLEO"

# Row 5
$ws.Range("B5").Value = "HuobiToken"
$ws.Range("C5").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "5.295"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "4HuobiTokenHT"

# Row 6
$ws.Range("B6").Value = "Cronos"
$ws.Range("C6").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.05794"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "5CronosCRO"

# Row 7
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.482"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "6KuCoinTokenKCS"

# Row 8
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.358"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "7GateTokenGT"

# Row 9
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8088"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "8MXTokenMX"

# Row 10
$ws.Range("B10").Value = "FTXToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8757"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "9FTXTokenFTT"

# Row 11
$ws.Range("B11").Value = "One"
$ws.Range("C11").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.01037"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "10OneONEBestin24h"

# Row 12
$ws.Range("B12").Value = "WazirX"
$ws.Range("C12").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1381"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "11WazirXWRX"

# Row 13
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07302"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "12MandalaExchangeTokenMDX"

# Row 14
$ws.Range("B14").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C14").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03066"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "13LiechtensteinCryptoassetsExchangeLCX"

# Row 15
$ws.Range("B15").Value = "BitrueCoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.03057"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "14BitrueCoinBTR"

# Row 16
$ws.Range("B16").Value = "BitMartToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.09334"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "15BitMartTokenBMX"

# Row 17
$ws.Range("B17").Value = "MCDex"
$ws.Range("C17").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.873"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "16MCDexMCB"

# Row 18
$ws.Range("B18").Value = "BitForexToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.001538"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "17BitForexTokenBF"

# Row 19
$ws.Range("B19").Value = "CoinExToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.04707"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "18CoinExTokenCET"

# Row 20
$ws.Range("B20").Value = "TigerCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.006087"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "19TigerCashTCH"

# Row 21
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001275"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "20BitKanKAN"

# Row 22
$ws.Range("B22").Value = "HotbitToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.004595"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "21HotbitTokenHTB"

# Row 23
$ws.Range("B23").Value = "NitroEx"
$ws.Range("C23").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.00008702"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "22NitroExNTX"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.140"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3211"

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0002344"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03773"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006401"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.004001"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007640"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005477"

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5901"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.006275"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "47BOLOBOLO"

